# Update the answer key table: replace the letter-grade values in the
# three data rows (rows 2-4 of the table, columns 3-7) with the new
# design's values. Row/column layout (1-indexed, Word convention):
#   Row 1: header (B\№, 1, 2, 3, 4, 5, 6)          -- unchanged
#   Row 2: "1", Г, В, Б, Г, В, Г  -> "1", Г, А, А, А, Г, Б
#   Row 3: "2", А, В, В, В, Г, А  -> "2", А, А, Б, Б, А, В
#   Row 4: "3", Б, А, Г, Г, Б, А  -> "3", Б, А, Г, Б, А, Б
$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(2,3).Range.Text = "А"
$t.Cell(2,4).Range.Text = "А"
$t.Cell(2,5).Range.Text = "А"
$t.Cell(2,6).Range.Text = "Г"
$t.Cell(2,7).Range.Text = "Б"

$t.Cell(3,3).Range.Text = "А"
$t.Cell(3,4).Range.Text = "Б"
$t.Cell(3,5).Range.Text = "Б"
$t.Cell(3,6).Range.Text = "А"
$t.Cell(3,7).Range.Text = "В"

$t.Cell(4,5).Range.Text = "Б"
$t.Cell(4,6).Range.Text = "А"
$t.Cell(4,7).Range.Text = "Б"
